# Approval MIS report functionality / loanEntry.missingfeields
# Adds a new "loandetails" worksheet (with sample row) at the end of the
# workbook, switches the active tab back to "partners" (with a new
# selection), and tidies up the selection left on the "logincrds" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "loandetails" worksheet after the last existing sheet
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "loandetails"

# Column widths roughly matching the authored layout.
$newSheet.Columns.Item(1).ColumnWidth = 17.36328125
$newSheet.Columns.Item(2).ColumnWidth = 26
$newSheet.Columns.Item(3).ColumnWidth = 25.90625
$newSheet.Columns.Item(4).ColumnWidth = 17.453125

# Header row - written in an order that reproduces the shared-string
# table layout of the authored workbook (new strings appended in the
# order: mobilenumber, partnercustomerid, testAutomationg001,
# 6000010000, groupid).
$newSheet.Range("A1").Value = "mobilenumber"
$newSheet.Range("B1").Value = "partnerloanid"
$newSheet.Range("C1").Value = "partnercustomerid"

# Sample data row.
$newSheet.Range("B2").Value = "testAutomationg001"
$newSheet.Range("C2").Value = "testAutomationg001"
$newSheet.Range("D2").Value = "testAutomationg001"
$newSheet.Range("A2").Value = "'6000010000"

$newSheet.Range("D1").Value = "groupid"

# Leave the new sheet's own selection on D4, as authored.
$newSheet.Range("D4").Select()

# ---------------------------------------------------------------------
# 2. Restore "partners" as the active tab and set its new selection.
# ---------------------------------------------------------------------
$partners = $wb.Worksheets.Item("partners")
$partners.Activate()
$partners.Range("J3").Select()
